$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 33336178
$ws.Range("I32").Value = 83335780
$ws.Range("J32").Value = 3115.6667
$ws.Range("K32").Value = 83335780
$ws.Range("L32").Value = 3115.6667
$ws.Range("M32").Value = -83335454
$ws.Range("N32").Value = -3767.6667
$ws.Range("H42").Value = 416.1111
$ws.Range("J42").Value = 439.83334
$ws.Range("L42").Value = 1319.50002
$ws.Range("N42").Value = -1779.50002
$ws.Range("H43").Value = 3086700.5
$ws.Range("I43").Value = 3857264.5
$ws.Range("K43").Value = 3857264.5
$ws.Range("M43").Value = -3857195.5
$ws.Range("H129").Value = 1820
$ws.Range("J129").Value = 3200
$ws.Range("L129").Value = 9600
$ws.Range("N129").Value = -19600
$ws.Range("H132").Value = 6836.2
$ws.Range("J132").Value = 13082.022
$ws.Range("L132").Value = 39246.06600000001
$ws.Range("N132").Value = -44306.06600000001
$ws.Range("H135").Value = 3221.8333
$ws.Range("I135").Value = 1952.3103
$ws.Range("K135").Value = 17570.7927
$ws.Range("M135").Value = -15035.7927
$ws.Range("H140").Value = 66921.42999999999
$ws.Range("J140").Value = 66408.336
$ws.Range("L140").Value = 66408.336
$ws.Range("N140").Value = -76768.336

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 755.5714
$ws.Range("I5").Value = 865.6667
$ws.Range("J5").Value = 95
$ws.Range("K5").Value = 865.6667
$ws.Range("L5").Value = 95
$ws.Range("M5").Value = -753.6667
$ws.Range("N5").Value = -319
$ws.Range("H74").Value = 3815.2974
$ws.Range("I74").Value = 1528.3
$ws.Range("J74").Value = 6505.8823
$ws.Range("K74").Value = 1528.3
$ws.Range("L74").Value = 6505.8823
$ws.Range("M74").Value = -654.3
$ws.Range("N74").Value = -8253.882300000001
$ws.Range("H77").Value = 3815.2974
$ws.Range("I77").Value = 1528.3
$ws.Range("J77").Value = 6505.8823
$ws.Range("K77").Value = 7641.5
$ws.Range("L77").Value = 32529.4115
$ws.Range("M77").Value = -3273.5
$ws.Range("N77").Value = -41265.4115
$ws.Range("H122").Value = 464808
$ws.Range("I122").Value = 921499.5
$ws.Range("K122").Value = 2764498.5
$ws.Range("M122").Value = -2762048.5
$ws.Range("H132").Value = 37833.844
$ws.Range("I132").Value = 44275.93
$ws.Range("J132").Value = 19796
$ws.Range("K132").Value = 132827.79
$ws.Range("L132").Value = 59388
$ws.Range("M132").Value = -130297.79
$ws.Range("N132").Value = -64448

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 755.5714
$ws.Range("I4").Value = 865.6667
$ws.Range("J4").Value = 95
$ws.Range("K4").Value = 865.6667
$ws.Range("L4").Value = 95
$ws.Range("M4").Value = -750.6667
$ws.Range("N4").Value = -325
$ws.Range("H105").Value = 142863280
$ws.Range("I105").Value = 166673010
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 166673010
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -166671263
$ws.Range("N105").Value = -8494
$ws.Range("H134").Value = 5484.3213
$ws.Range("I134").Value = 2381.182
$ws.Range("J134").Value = 7492.2354
$ws.Range("K134").Value = 7143.545999999999
$ws.Range("L134").Value = 22476.7062
$ws.Range("M134").Value = -4608.545999999999
$ws.Range("N134").Value = -27546.7062

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 19950
$ws.Range("I45").Value = 19950
$ws.Range("K45").Value = 19950
$ws.Range("M45").Value = -19357
$ws.Range("H132").Value = 17861118
$ws.Range("I132").Value = 18870134
$ws.Range("K132").Value = 56610402
$ws.Range("M132").Value = -56607872
$ws.Range("H134").Value = 2026.5238
$ws.Range("I134").Value = 1952.85
$ws.Range("K134").Value = 5858.549999999999
$ws.Range("M134").Value = -3323.549999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 590.5217
$ws.Range("I5").Value = 572.2857
$ws.Range("J5").Value = 618.8889
$ws.Range("K5").Value = 1716.8571
$ws.Range("L5").Value = 1856.6667
$ws.Range("M5").Value = -1604.8571
$ws.Range("N5").Value = -2080.6667
$ws.Range("H34").Value = 645825.4399999999
$ws.Range("J34").Value = 1595.2667
$ws.Range("L34").Value = 4785.800099999999
$ws.Range("N34").Value = -4953.800099999999
$ws.Range("H39").Value = 1592.5555
$ws.Range("J39").Value = 998
$ws.Range("L39").Value = 2994
$ws.Range("N39").Value = -3582
$ws.Range("H55").Value = 1500
$ws.Range("J55").Value = 1500
$ws.Range("L55").Value = 4500
$ws.Range("N55").Value = -4854
$ws.Range("H107").Value = 354669.2
$ws.Range("I107").Value = 1189.6
$ws.Range("K107").Value = 3568.8
$ws.Range("M107").Value = -1648.8
$ws.Range("H134").Value = 3900.05
$ws.Range("I134").Value = 3529.4707
$ws.Range("K134").Value = 10588.4121
$ws.Range("M134").Value = -5518.4121
$ws.Range("H135").Value = 590.5217
$ws.Range("I135").Value = 572.2857
$ws.Range("J135").Value = 618.8889
$ws.Range("K135").Value = 5150.571300000001
$ws.Range("L135").Value = 5570.0001
$ws.Range("M135").Value = -2615.571300000001
$ws.Range("N135").Value = -10640.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 44367.332
$ws.Range("J93").Value = 44367.332
$ws.Range("L93").Value = 44367.332
$ws.Range("N93").Value = -48111.332
$ws.Range("H113").Value = 16725194
$ws.Range("I113").Value = 18583272
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 18583272
$ws.Range("L113").Value = 2499.5
$ws.Range("M113").Value = -18581102
$ws.Range("N113").Value = -6839.5
$ws.Range("H122").Value = 505828.12
$ws.Range("I122").Value = 850861.9
$ws.Range("K122").Value = 2552585.7
$ws.Range("M122").Value = -2550135.7
$ws.Range("H132").Value = 4573.3716
$ws.Range("I132").Value = 4501.759
$ws.Range("K132").Value = 13505.277
$ws.Range("M132").Value = -10975.277

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4216.03
$ws.Range("I132").Value = 4216.03
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12648.09
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10118.09
$ws.Range("N132").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H62").Value = 7000
$ws.Range("J62").Value = 7000
$ws.Range("L62").Value = 7000
$ws.Range("N62").Value = -8248
$ws.Range("H65").Value = 7000
$ws.Range("J65").Value = 7000
$ws.Range("L65").Value = 35000
$ws.Range("N65").Value = -41240
$ws.Range("H126").Value = 4663.625
$ws.Range("I126").Value = 4061.4
$ws.Range("K126").Value = 12184.2
$ws.Range("M126").Value = -9714.200000000001
$ws.Range("H132").Value = 2424191.2
$ws.Range("I132").Value = 3714294.5
$ws.Range("J132").Value = 5248.125
$ws.Range("K132").Value = 11142883.5
$ws.Range("L132").Value = 15744.375
$ws.Range("M132").Value = -11140353.5
$ws.Range("N132").Value = -20804.375
$ws.Range("H136").Value = 8561.354499999999
$ws.Range("I136").Value = 3583.25
$ws.Range("J136").Value = 9556.975
$ws.Range("K136").Value = 10749.75
$ws.Range("L136").Value = 28670.925
$ws.Range("M136").Value = -8199.75
$ws.Range("N136").Value = -33770.925
